$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 2.341355666666667
$ws.Range("N2").Value = 7.024067000000001
$ws.Range("O2").Value = 0.03973512964576821
$ws.Range("P2").Value = 0.0397351296457682
$ws.Range("Q2").Value = 10.01665123405278
$ws.Range("R2").Value = 90.149861106475
$ws.Range("S2").Value = 0.03839695019198483
$ws.Range("T2").Value = 0.03839695019198482

# Row 3
$ws.Range("O3").Value = 0.5779093692199981
$ws.Range("P3").Value = 0.5779093692199981
$ws.Range("S3").Value = 0.5584468318900997
$ws.Range("T3").Value = 0.5584468318900997

# Row 4
$ws.Range("O4").Value = 0.3823555011342337
$ws.Range("P4").Value = 0.3823555011342337
$ws.Range("S4").Value = 0.3694787273519347
$ws.Range("T4").Value = 0.3694787273519346

# Row 5
$ws.Range("M5").Value = 2.341355666666667
$ws.Range("N5").Value = 7.024067000000001
$ws.Range("O5").Value = 0.03973512964576821
$ws.Range("P5").Value = 0.0397351296457682
$ws.Range("Q5").Value = 0.3490922276405556
$ws.Range("R5").Value = 3.141830048765
$ws.Range("S5").Value = 0.001338179453783388
$ws.Range("T5").Value = 0.001338179453783387

# Row 6
$ws.Range("O6").Value = 0.5779093692199981
$ws.Range("P6").Value = 0.5779093692199981
$ws.Range("Q6").Value = 5.077211799077222
$ws.Range("S6").Value = 0.01946253732989847
$ws.Range("T6").Value = 0.01946253732989847

# Row 7
$ws.Range("O7").Value = 0.3823555011342337
$ws.Range("P7").Value = 0.3823555011342337
$ws.Range("S7").Value = 0.01287677378229906
$ws.Range("T7").Value = 0.01287677378229906
